$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '91.686.25'
$ws.Range("E2").Value = '  +1.25%  '
$ws.Range("D3").Value = '3.129.58'
$ws.Range("E3").Value = '  -0.42%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '241.60'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.53%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '626.35'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.79%  '
$ws.Range("E7").Value = '  +7.13%  '
$ws.Range("E8").Value = '  +4.30%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.999'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.11%  '
$ws.Range("D10").Value = '3.128.74'
$ws.Range("E10").Value = '  -0.26%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.760'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +5.36%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.204'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.82%  '
$ws.Range("E13").Value = '  +3.23%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.83'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.41%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.53'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.16%  '
$ws.Range("D16").Value = '91.274.46'
$ws.Range("E16").Value = '  +1.19%  '
$ws.Range("D17").Value = '3.706.71'
$ws.Range("E17").Value = '  -0.11%  '
$ws.Range("D18").Value = '3.151.96'
$ws.Range("E18").Value = '  +1.05%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.81'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.44%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.71'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.32%  '
$ws.Range("E21").Value = '  -0.37%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.89'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.50%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '452.18'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.30%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.21'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.59%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.95'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.53%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '93.26'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.87%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.09'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.20%  '
$ws.Range("E29").Value = '  +0.11%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.182'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +13.29%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.234'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +16.21%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.119'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +38.35%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '9.26'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -7.26%  '
$ws.Range("E34").Value = '  +36.36%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.164'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +9.05%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '26.93'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.49%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '7.57'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +5.86%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.18'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +26.07%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '500.09'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.66%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.94'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.04%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.64'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -6.50%  '
$ws.Range("E42").Value = '  -0.68%  '
$ws.Range("E43").Value = '  +0.00%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '22.17'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.17%  '
$ws.Range("E46").Value = '  +0.16%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '157.61'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +5.14%  '
$ws.Range("E48").Value = '  -0.38%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '4.58'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.46%  '
$ws.Range("E50").Value = '  +0.22%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '44.88'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.58%  '